$d = $word.ActiveDocument
$rng = $d.Range(1090, 1095)
Write-Output ("Range Text: [" + $rng.Text + "]")
$rng.InsertParagraphAfter()

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Note:`r") {
        Write-Output ("Isolated para: [" + $p.Range.Text + "] start=" + $p.Range.Start + " end=" + $p.Range.End)
        $p.Range.Font.BoldBi = $true
        # Now delete the paragraph mark to merge back with next paragraph
        $markRng = $d.Range($p.Range.End - 1, $p.Range.End)
        Write-Output ("Mark range text code: [" + $markRng.Text + "]")
        $markRng.Delete()
        break
    }
}
